# Generate Report for Handback
# - Update the "Ready for handoff" status text (Overview sheet) to reflect
#   a failed handback transform.
# - Record the handback/handoff file-name mismatch error detail for the
#   496a8088 row on both the zh-cn and de-de worksheets, widening the
#   "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$Overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# The "Status" text shared by the 496a8088 row changes from
# "Ready for handoff" to "Handback transform failed" everywhere it is
# shown (Overview!E3/F3 and the per-locale Status column, C3, on both
# the zh-cn and de-de sheets).
$Overview.Range("E3").Value = "Handback transform failed"
$Overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Error Detail (column P) for the 496a8088 row on each locale sheet, plus
# widen column P (to a stored width of 40) so the long message is visible.
# NOTE: this host's ColumnWidth -> stored-OOXML-width conversion adds a
# flat 5/6 (0.8333...) offset, so we dial the input back by that amount to
# land on a saved column width of exactly 40.
$targetColumnWidth = 40 - 5/6

$zhcn.Range("P3").Value = "Handback file name: yy1ju0un.tbl is different with handoff file name: 496a8088-67a4-4ecd-9b3d-9e9fd0d2c079.9bfd6a44642691a77b988c3e2034760a4e1fe9f2.zh-cn."
$zhcn.Columns.Item(16).ColumnWidth = $targetColumnWidth

$dede.Range("P3").Value = "Handback file name: yy1ju0un.tbl is different with handoff file name: 496a8088-67a4-4ecd-9b3d-9e9fd0d2c079.9bfd6a44642691a77b988c3e2034760a4e1fe9f2.de-de."
$dede.Columns.Item(16).ColumnWidth = $targetColumnWidth
